{"js": "const pairs = [\n  [\"2024-01-23 Tuesday\", \"2024-01-24 Wednesday\"],\n  [\"47+3=50\", \"19+12=31\"],\n  [\"94-18=76\", \"99-78=21\"],\n  [\"9+14=23\", \"90-52=38\"],\n  [\"2+38=40\", \"40-3=37\"],\n  [\"21-13=8\", \"20+61=81\"],\n  [\"56-54=2\", \"98-39=59\"],\n  [\"12+83=95\", \"53-2=51\"],\n  [\"28-11=17\", \"97-18=79\"],\n  [\"84-45=39\", \"60-22=38\"],\n  [\"41+55=96\", \"6+56=62\"],\n  [\"7+29=36\", \"26+36=62\"],\n  [\"65-60=5\", \"83-4=79\"],\n  [\"12+76=88\", \"6-6=0\"],\n  [\"38+31=69\", \"90+9=99\"],\n  [\"17+70=87\", \"83-5=78\"],\n  [\"27+36=63\", \"57+16=73\"],\n  [\"7+55=62\", \"3+67=70\"],\n  [\"70+21=91\", \"67+8=75\"],\n  [\"54-36=18\", \"61+36=97\"],\n  [\"82-23=59\", \"37+27=64\"],\n  [\"73-10=63\", \"31+21=52\"],\n  [\"69+7=76\", \"66+27=93\"],\n  [\"6+12=18\", \"17+52=69\"],\n  [\"9+89=98\", \"66-65=1\"],\n  [\"91-61=30\", \"29+36=65\"],\n  [\"79-57=22\", \"13+5=18\"],\n  [\"27+51=78\", \"90-87=3\"],\n  [\"72-19=53\", \"53-16=37\"],\n  [\"32-13=19\", \"96-86=10\"],\n  [\"86-18=68\", \"72-5=67\"],\n  [\"22+6=28\", \"67-10=57\"],\n  [\"32-0=32\", \"52+16=68\"],\n  [\"62+4=66\", \"45-21=24\"],\n  [\"15+62=77\", \"93-44=49\"],\n  [\"24-6=18\", \"8+43=51\"],\n  [\"48-18=30\", \"53+21=74\"],\n  [\"24+20=44\", \"78+21=99\"],\n  [\"5+29=34\", \"67-46=21\"],\n  [\"17+13=30\", \"36-8=28\"],\n  [\"8+23=31\", \"93-27=66\"],\n  [\"69-33=36\", \"73+5=78\"],\n  [\"88-73=15\", \"77-36=41\"],\n  [\"20+14=34\", \"17+73=90\"],\n  [\"95-85=10\", \"11+72=83\"],\n  [\"93-30=63\", \"29+39=68\"],\n  [\"74+9=83\", \"25-20=5\"],\n  [\"15+57=72\", \"71-20=51\"],\n  [\"73-0=73\", \"39-20=19\"],\n  [\"39-4=35\", \"55-17=38\"],\n  [\"46+38=84\", \"84-70=14\"],\n  [\"29+0=29\", \"30-14=16\"],\n  [\"98-16=82\", \"71-27=44\"],\n  [\"76+13=89\", \"17+81=98\"],\n  [\"0+58=58\", \"81+10=91\"],\n  [\"17+10=27\", \"11+82=93\"],\n  [\"13+33=46\", \"78+2=80\"],\n  [\"59+38=97\", \"13+51=64\"],\n  [\"26+12=38\", \"45+15=60\"],\n  [\"50+4=54\", \"6+46=52\"],\n  [\"77-6=71\", \"61+24=85\"],\n  [\"36+12=48\", \"4+35=39\"],\n  [\"72-15=57\", \"49+5=54\"],\n  [\"19+69=88\", \"74-11=63\"],\n  [\"86-43=43\", \"64-32=32\"],\n  [\"42+51=93\", \"14+45=59\"],\n  [\"67+29=96\", \"59+18=77\"],\n  [\"74-43=31\", \"78-61=17\"],\n  [\"78-58=20\", \"66-0=66\"],\n  [\"12+82=94\", \"57-45=12\"],\n  [\"83-59=24\", \"59-58=1\"],\n  [\"94-41=53\", \"91-61=30\"],\n  [\"80-16=64\", \"14+55=69\"],\n  [\"29+66=95\", \"85-33=52\"],\n  [\"26+8=34\", \"0+60=60\"],\n  [\"47-25=22\", \"89-15=74\"],\n  [\"28-20=8\", \"3+29=32\"],\n  [\"60+13=73\", \"93-5=88\"],\n  [\"20+68=88\", \"40+3=43\"],\n  [\"93-0=93\", \"60+3=63\"],\n  [\"9+18=27\", \"39-11=28\"],\n  [\"84-71=13\", \"84-25=59\"],\n  [\"75+21=96\", \"83-0=83\"],\n  [\"49-19=30\", \"38-6=32\"],\n  [\"1+16=17\", \"15+44=59\"],\n  [\"94-75=19\", \"97-54=43\"],\n  [\"78-6=72\", \"31+7=38\"],\n  [\"46+25=71\", \"30+5=35\"],\n  [\"33+66=99\", \"51-40=11\"],\n  [\"91-50=41\", \"49+11=60\"],\n  [\"14+46=60\", \"2+91=93\"],\n  [\"68-50=18\", \"65-26=39\"],\n  [\"62-20=42\", \"73+19=92\"],\n  [\"54-6=48\", \"20-11=9\"],\n  [\"91-34=57\", \"33+40=73\"],\n  [\"52+32=84\", \"67+10=77\"],\n  [\"5+77=82\", \"22+33=55\"],\n  [\"91+0=91\", \"98-30=68\"],\n  [\"92-11=81\", \"90-12=78\"],\n  [\"14+23=37\", \"80-36=44\"],\n  [\"17-8=9\", \"49+3=52\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2024-01-23 Tuesday\", \"2024-01-24 Wednesday\")\n  ,@(\"47+3=50\", \"19+12=31\")\n  ,@(\"94-18=76\", \"99-78=21\")\n  ,@(\"9+14=23\", \"90-52=38\")\n  ,@(\"2+38=40\", \"40-3=37\")\n  ,@(\"21-13=8\", \"20+61=81\")\n  ,@(\"56-54=2\", \"98-39=59\")\n  ,@(\"12+83=95\", \"53-2=51\")\n  ,@(\"28-11=17\", \"97-18=79\")\n  ,@(\"84-45=39\", \"60-22=38\")\n  ,@(\"41+55=96\", \"6+56=62\")\n  ,@(\"7+29=36\", \"26+36=62\")\n  ,@(\"65-60=5\", \"83-4=79\")\n  ,@(\"12+76=88\", \"6-6=0\")\n  ,@(\"38+31=69\", \"90+9=99\")\n  ,@(\"17+70=87\", \"83-5=78\")\n  ,@(\"27+36=63\", \"57+16=73\")\n  ,@(\"7+55=62\", \"3+67=70\")\n  ,@(\"70+21=91\", \"67+8=75\")\n  ,@(\"54-36=18\", \"61+36=97\")\n  ,@(\"82-23=59\", \"37+27=64\")\n  ,@(\"73-10=63\", \"31+21=52\")\n  ,@(\"69+7=76\", \"66+27=93\")\n  ,@(\"6+12=18\", \"17+52=69\")\n  ,@(\"9+89=98\", \"66-65=1\")\n  ,@(\"91-61=30\", \"29+36=65\")\n  ,@(\"79-57=22\", \"13+5=18\")\n  ,@(\"27+51=78\", \"90-87=3\")\n  ,@(\"72-19=53\", \"53-16=37\")\n  ,@(\"32-13=19\", \"96-86=10\")\n  ,@(\"86-18=68\", \"72-5=67\")\n  ,@(\"22+6=28\", \"67-10=57\")\n  ,@(\"32-0=32\", \"52+16=68\")\n  ,@(\"62+4=66\", \"45-21=24\")\n  ,@(\"15+62=77\", \"93-44=49\")\n  ,@(\"24-6=18\", \"8+43=51\")\n  ,@(\"48-18=30\", \"53+21=74\")\n  ,@(\"24+20=44\", \"78+21=99\")\n  ,@(\"5+29=34\", \"67-46=21\")\n  ,@(\"17+13=30\", \"36-8=28\")\n  ,@(\"8+23=31\", \"93-27=66\")\n  ,@(\"69-33=36\", \"73+5=78\")\n  ,@(\"88-73=15\", \"77-36=41\")\n  ,@(\"20+14=34\", \"17+73=90\")\n  ,@(\"95-85=10\", \"11+72=83\")\n  ,@(\"93-30=63\", \"29+39=68\")\n  ,@(\"74+9=83\", \"25-20=5\")\n  ,@(\"15+57=72\", \"71-20=51\")\n  ,@(\"73-0=73\", \"39-20=19\")\n  ,@(\"39-4=35\", \"55-17=38\")\n  ,@(\"46+38=84\", \"84-70=14\")\n  ,@(\"29+0=29\", \"30-14=16\")\n  ,@(\"98-16=82\", \"71-27=44\")\n  ,@(\"76+13=89\", \"17+81=98\")\n  ,@(\"0+58=58\", \"81+10=91\")\n  ,@(\"17+10=27\", \"11+82=93\")\n  ,@(\"13+33=46\", \"78+2=80\")\n  ,@(\"59+38=97\", \"13+51=64\")\n  ,@(\"26+12=38\", \"45+15=60\")\n  ,@(\"50+4=54\", \"6+46=52\")\n  ,@(\"77-6=71\", \"61+24=85\")\n  ,@(\"36+12=48\", \"4+35=39\")\n  ,@(\"72-15=57\", \"49+5=54\")\n  ,@(\"19+69=88\", \"74-11=63\")\n  ,@(\"86-43=43\", \"64-32=32\")\n  ,@(\"42+51=93\", \"14+45=59\")\n  ,@(\"67+29=96\", \"59+18=77\")\n  ,@(\"74-43=31\", \"78-61=17\")\n  ,@(\"78-58=20\", \"66-0=66\")\n  ,@(\"12+82=94\", \"57-45=12\")\n  ,@(\"83-59=24\", \"59-58=1\")\n  ,@(\"94-41=53\", \"91-61=30\")\n  ,@(\"80-16=64\", \"14+55=69\")\n  ,@(\"29+66=95\", \"85-33=52\")\n  ,@(\"26+8=34\", \"0+60=60\")\n  ,@(\"47-25=22\", \"89-15=74\")\n  ,@(\"28-20=8\", \"3+29=32\")\n  ,@(\"60+13=73\", \"93-5=88\")\n  ,@(\"20+68=88\", \"40+3=43\")\n  ,@(\"93-0=93\", \"60+3=63\")\n  ,@(\"9+18=27\", \"39-11=28\")\n  ,@(\"84-71=13\", \"84-25=59\")\n  ,@(\"75+21=96\", \"83-0=83\")\n  ,@(\"49-19=30\", \"38-6=32\")\n  ,@(\"1+16=17\", \"15+44=59\")\n  ,@(\"94-75=19\", \"97-54=43\")\n  ,@(\"78-6=72\", \"31+7=38\")\n  ,@(\"46+25=71\", \"30+5=35\")\n  ,@(\"33+66=99\", \"51-40=11\")\n  ,@(\"91-50=41\", \"49+11=60\")\n  ,@(\"14+46=60\", \"2+91=93\")\n  ,@(\"68-50=18\", \"65-26=39\")\n  ,@(\"62-20=42\", \"73+19=92\")\n  ,@(\"54-6=48\", \"20-11=9\")\n  ,@(\"91-34=57\", \"33+40=73\")\n  ,@(\"52+32=84\", \"67+10=77\")\n  ,@(\"5+77=82\", \"22+33=55\")\n  ,@(\"91+0=91\", \"98-30=68\")\n  ,@(\"92-11=81\", \"90-12=78\")\n  ,@(\"14+23=37\", \"80-36=44\")\n  ,@(\"17-8=9\", \"49+3=52\")\n)\n\nforeach ($p in $pairs) {\n  $old = $p[0]\n  $new = $p[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    throw \"No match found for: $old\"\n  }\n}"}
